$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.056.83"
$ws.Range("E2").Value = "  +15.52%  "
$ws.Range("D3").Value = "1.668.65"
$ws.Range("E3").Value = "  +10.01%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'307.30"
$ws.Range("E5").Value = "  +9.81%  "
$ws.Range("D6").Value = "'0.9957"
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("D7").Value = "'0.3707"
$ws.Range("E7").Value = "  +3.37%  "
$ws.Range("D8").Value = "'0.3432"
$ws.Range("E8").Value = "  +9.59%  "
$ws.Range("D9").Value = "'48.06"
$ws.Range("E9").Value = "  +20.22%  "
$ws.Range("D10").Value = "'1.170"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").Value = "'0.07236"
$ws.Range("E11").Value = "  +7.49%  "
$ws.Range("D12").Value = "'0.9994"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "'20.43"
$ws.Range("E13").Value = "  +9.70%  "
$ws.Range("D14").Value = "'6.019"
$ws.Range("E14").Value = "  +7.62%  "
$ws.Range("D15").Value = "'6.728"
$ws.Range("E15").Value = "  +7.45%  "
$ws.Range("D16").Value = "1.667.88"
$ws.Range("E16").Value = "  +10.50%  "
$ws.Range("D17").Value = "'0.00001098"
$ws.Range("E17").Value = "  +6.62%  "
$ws.Range("D18").Value = "'0.9958"
$ws.Range("E18").Value = "  +3.80%  "
$ws.Range("D19").Value = "'0.06701"
$ws.Range("E19").Value = "  +10.96%  "
$ws.Range("D20").Value = "'81.67"
$ws.Range("E20").Value = "  +16.09%  "
$ws.Range("D21").Value = "'16.40"
$ws.Range("E21").Value = "  +10.30%  "
$ws.Range("D22").Value = "'6.125"
$ws.Range("E22").Value = "  +9.28%  "
$ws.Range("D23").Value = "'11.96"
$ws.Range("E23").Value = "  +5.47%  "
$ws.Range("D24").Value = "24.041.34"
$ws.Range("E24").Value = "  +15.33%  "
$ws.Range("D25").Value = "'2.400"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").Value = "'3.385"
$ws.Range("E26").Value = "  -7.81%  "
$ws.Range("D27").Value = "'2.659"
$ws.Range("E27").Value = "  +23.53%  "
$ws.Range("D28").Value = "'152.07"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").Value = "'19.52"
$ws.Range("E29").Value = "  +11.73%  "
$ws.Range("D30").Value = "1.849.71"
$ws.Range("E30").Value = "  +10.63%  "
$ws.Range("D31").Value = "'127.00"
$ws.Range("E31").Value = "  +9.42%  "
$ws.Range("D32").Value = "'6.351"
$ws.Range("E32").Value = "  +24.59%  "
$ws.Range("D33").Value = "'4.070"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").Value = "'0.9787"
$ws.Range("E34").Value = "  +16.88%  "
$ws.Range("D35").Value = "'1.731"
$ws.Range("E35").Value = "  +17.77%  "
$ws.Range("D36").Value = "'0.08373"
$ws.Range("E36").Value = "  +4.65%  "
$ws.Range("D37").Value = "'12.27"
$ws.Range("E37").Value = "  +16.40%  "
$ws.Range("D38").Value = "'8.916"
$ws.Range("E38").Value = "  +18.36%  "
$ws.Range("D39").Value = "'5.311"
$ws.Range("E39").Value = "  +9.51%  "
$ws.Range("D40").Value = "'0.06346"
$ws.Range("E40").Value = "  +9.49%  "
$ws.Range("D41").Value = "'1.289"
$ws.Range("E41").Value = "  +6.35%  "
$ws.Range("D42").Value = "'0.02314"
$ws.Range("E42").Value = "  +11.54%  "
$ws.Range("D43").Value = "'0.2074"
$ws.Range("E43").Value = "  +10.12%  "
$ws.Range("D44").Value = "'0.6081"
$ws.Range("E44").Value = "  +14.15%  "
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").Value = "'3.802"
$ws.Range("E46").Value = "  +7.33%  "
$ws.Range("D47").Value = "'13.19"
$ws.Range("E47").Value = "  +7.45%  "
$ws.Range("D48").Value = "'0.5927"
$ws.Range("E48").Value = "  +11.22%  "
$ws.Range("D49").Value = "'127.03"
$ws.Range("E49").Value = "  +5.45%  "
$ws.Range("D50").Value = "'2.001"
$ws.Range("E50").Value = "  +7.88%  "
$ws.Range("D51").Value = "'0.07109"
$ws.Range("E51").Value = "  +9.31%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
